# Row 7 and Row 8 of the "Artfynd" sheet swap their per-record data.
# A handful of columns happen to hold identical values in both rows
# (AA, AD, AE, AG, AT, AW, AY, D, S, T, U, V, W, Y) so they are left
# untouched - swapping identical values would be a no-op anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7 gets what used to be Row 8's record ----
$ws.Range("A7").Value = 131106436
$ws.Range("B7").Value = 5493
$ws.Range("E7").Value = 101410
$ws.Range("F7").Value = "Reliktbock"
$ws.Range("G7").Value = "Nothorhina muricata"
$ws.Range("H7").Value = "(Dalman, 1817)"
$ws.Range("I7").Value = "'2"
$ws.Range("J7").Value = "ex."
$ws.Range("P7").Value = "Svartmyran, Mpd"
$ws.Range("Q7").Value = 616762
$ws.Range("R7").Value = 6934714
$ws.Range("X7").Value = "2025_0743"
$ws.Range("Z7").Value = "11:39"
$ws.Range("AB7").Value = "11:39"
$ws.Range("AC7").Value = "Två kläckhål"
$ws.Range("AX7").Value = "David Isaksson"

# ---- Row 8 gets what used to be Row 7's record ----
$ws.Range("A8").Value = 131108352
$ws.Range("B8").Value = 80214
$ws.Range("E8").Value = 388
$ws.Range("F8").Value = "Stiftgelélav"
$ws.Range("G8").Value = "Collema furfuraceum"
$ws.Range("H8").Value = "(Arnold) Du Rietz"
$ws.Range("I8").Value = "'1"
$ws.Range("J8").Value = "bålar"
$ws.Range("P8").Value = "S Svartmyran, Mpd"
$ws.Range("Q8").Value = 616863
$ws.Range("R8").Value = 6934788
$ws.Range("X8").Value = "2025_0758"
$ws.Range("Z8").Value = "14:47"
$ws.Range("AB8").Value = "14:47"
$ws.Range("AC8").Value = "Asp"
$ws.Range("AX8").Value = "Måns Svensson"

# The leading apostrophes on I7/I8 force Excel to store the numeric-
# looking "1"/"2" values as text (matching the source inlineStr type
# instead of being auto-coerced to a number). That quote-prefix also
# marks the cell style dirty, so reset it back to Normal/General
# afterwards to keep the style index identical to the original file.
$ws.Range("I7:I8").Style = "Normal"
